$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot original row data (columns B:AD) for every affected row
$snapshot = @{}
$snapshot[14] = $ws.Range("B14:AD14").Value()
$snapshot[15] = $ws.Range("B15:AD15").Value()
$snapshot[20] = $ws.Range("B20:AD20").Value()
$snapshot[22] = $ws.Range("B22:AD22").Value()
$snapshot[24] = $ws.Range("B24:AD24").Value()
$snapshot[25] = $ws.Range("B25:AD25").Value()
$snapshot[26] = $ws.Range("B26:AD26").Value()
$snapshot[27] = $ws.Range("B27:AD27").Value()
$snapshot[30] = $ws.Range("B30:AD30").Value()
$snapshot[31] = $ws.Range("B31:AD31").Value()
$snapshot[32] = $ws.Range("B32:AD32").Value()
$snapshot[33] = $ws.Range("B33:AD33").Value()
$snapshot[34] = $ws.Range("B34:AD34").Value()
$snapshot[35] = $ws.Range("B35:AD35").Value()
$snapshot[36] = $ws.Range("B36:AD36").Value()
$snapshot[37] = $ws.Range("B37:AD37").Value()
$snapshot[77] = $ws.Range("B77:AD77").Value()
$snapshot[78] = $ws.Range("B78:AD78").Value()
$snapshot[79] = $ws.Range("B79:AD79").Value()
$snapshot[80] = $ws.Range("B80:AD80").Value()
$snapshot[85] = $ws.Range("B85:AD85").Value()
$snapshot[87] = $ws.Range("B87:AD87").Value()
$snapshot[90] = $ws.Range("B90:AD90").Value()
$snapshot[91] = $ws.Range("B91:AD91").Value()
$snapshot[92] = $ws.Range("B92:AD92").Value()
$snapshot[95] = $ws.Range("B95:AD95").Value()
$snapshot[96] = $ws.Range("B96:AD96").Value()
$snapshot[141] = $ws.Range("B141:AD141").Value()
$snapshot[144] = $ws.Range("B144:AD144").Value()
$snapshot[145] = $ws.Range("B145:AD145").Value()
$snapshot[148] = $ws.Range("B148:AD148").Value()
$snapshot[155] = $ws.Range("B155:AD155").Value()
$snapshot[156] = $ws.Range("B156:AD156").Value()
$snapshot[157] = $ws.Range("B157:AD157").Value()
$snapshot[174] = $ws.Range("B174:AD174").Value()
$snapshot[175] = $ws.Range("B175:AD175").Value()
$snapshot[176] = $ws.Range("B176:AD176").Value()

# Write back according to the permutation (row <- source row)
$ws.Range("B14:AD14").Value = $snapshot[15]
$ws.Range("B15:AD15").Value = $snapshot[14]
$ws.Range("B20:AD20").Value = $snapshot[22]
$ws.Range("B22:AD22").Value = $snapshot[20]
$ws.Range("B24:AD24").Value = $snapshot[25]
$ws.Range("B25:AD25").Value = $snapshot[24]
$ws.Range("B26:AD26").Value = $snapshot[27]
$ws.Range("B27:AD27").Value = $snapshot[26]
$ws.Range("B30:AD30").Value = $snapshot[32]
$ws.Range("B31:AD31").Value = $snapshot[30]
$ws.Range("B32:AD32").Value = $snapshot[31]
$ws.Range("B33:AD33").Value = $snapshot[37]
$ws.Range("B34:AD34").Value = $snapshot[33]
$ws.Range("B35:AD35").Value = $snapshot[34]
$ws.Range("B36:AD36").Value = $snapshot[35]
$ws.Range("B37:AD37").Value = $snapshot[36]
$ws.Range("B77:AD77").Value = $snapshot[78]
$ws.Range("B78:AD78").Value = $snapshot[77]
$ws.Range("B79:AD79").Value = $snapshot[80]
$ws.Range("B80:AD80").Value = $snapshot[79]
$ws.Range("B85:AD85").Value = $snapshot[87]
$ws.Range("B87:AD87").Value = $snapshot[85]
$ws.Range("B90:AD90").Value = $snapshot[91]
$ws.Range("B91:AD91").Value = $snapshot[92]
$ws.Range("B92:AD92").Value = $snapshot[90]
$ws.Range("B95:AD95").Value = $snapshot[96]
$ws.Range("B96:AD96").Value = $snapshot[95]
$ws.Range("B141:AD141").Value = $snapshot[144]
$ws.Range("B144:AD144").Value = $snapshot[141]
$ws.Range("B145:AD145").Value = $snapshot[148]
$ws.Range("B148:AD148").Value = $snapshot[145]
$ws.Range("B155:AD155").Value = $snapshot[157]
$ws.Range("B156:AD156").Value = $snapshot[155]
$ws.Range("B157:AD157").Value = $snapshot[156]
$ws.Range("B174:AD174").Value = $snapshot[176]
$ws.Range("B175:AD175").Value = $snapshot[174]
$ws.Range("B176:AD176").Value = $snapshot[175]
